# Insert a new data row at row 454 (pushes existing rows 454-475 down to 455-476)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 454, shifting rows 454..475 down to 455..476
$ws.Rows.Item(454).Insert()

# Populate the newly inserted row 454 with the new record's data
$ws.Cells.Item(454, 1).Value = 4
$ws.Cells.Item(454, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(454, 3).Value = "Los Lagos"
$ws.Cells.Item(454, 4).Value = 45267
$ws.Cells.Item(454, 5).Value = 10
$ws.Cells.Item(454, 6).Value = 100112044
$ws.Cells.Item(454, 7).Value = "Perejil"
$ws.Cells.Item(454, 8).Value = "Sin especificar"
$ws.Cells.Item(454, 9).Value = "Primera"
$ws.Cells.Item(454, 10).Value = 50
$ws.Cells.Item(454, 11).Value = 8000
$ws.Cells.Item(454, 12).Value = 8000
$ws.Cells.Item(454, 13).Value = 8000
$ws.Cells.Item(454, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(454, 15).Value = "Región Metropolitana"
$ws.Cells.Item(454, 16).Value = 2667
$ws.Cells.Item(454, 17).Value = 3
$ws.Cells.Item(454, 18).Value = "Hortaliza"

Write-Host "Row 454 inserted and populated."
